$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as literal text
# (mirrors typing into a cell pre-formatted as Text), then restore the
# cell's original (unstyled) appearance by copying the style from a
# neighboring cell in column B that has no explicit style applied.
function Set-TextValue($cellRef, $rowNum, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $ws.Range("B$rowNum").Style
}

$ws.Range("D2").Value = "51.783.85"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "2.941.89"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue "D5" 5 "353.18"
$ws.Range("E5").Value = "  -1.45%  "

Set-TextValue "D6" 6 "105.52"
$ws.Range("E6").Value = "  -4.58%  "

Set-TextValue "D7" 7 "0.549"
$ws.Range("E7").Value = "  -3.62%  "

$ws.Range("E8").Value = "  +0.01%  "

Set-TextValue "D9" 9 "0.597"
$ws.Range("E9").Value = "  -5.65%  "

Set-TextValue "D10" 10 "37.43"
$ws.Range("E10").Value = "  -5.04%  "

$ws.Range("E11").Value = "  +1.99%  "

Set-TextValue "D12" 12 "0.0846"
$ws.Range("E12").Value = "  -3.75%  "

Set-TextValue "D13" 13 "18.78"
$ws.Range("E13").Value = "  -4.57%  "

$ws.Range("D14").Value = "3.402.20"
$ws.Range("E14").Value = "  +0.33%  "

Set-TextValue "D15" 15 "7.44"
$ws.Range("E15").Value = "  -5.92%  "

$ws.Range("D16").Value = "2.932.11"
$ws.Range("E16").Value = "  +0.50%  "

Set-TextValue "D17" 17 "0.981"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "51.712.18"
$ws.Range("E18").Value = "  -0.50%  "

Set-TextValue "D19" 19 "3.33"
$ws.Range("E19").Value = "  -1.81%  "

Set-TextValue "D20" 20 "7.30"
$ws.Range("E20").Value = "  -4.13%  "

Set-TextValue "D21" 21 "13.23"
$ws.Range("E21").Value = "  -6.08%  "

$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  -3.09%  "

Set-TextValue "D23" 23 "68.83"
$ws.Range("E23").Value = "  -3.26%  "

Set-TextValue "D24" 24 "264.88"
$ws.Range("E24").Value = "  -2.06%  "

Set-TextValue "D25" 25 "2.68"
$ws.Range("E25").Value = "  -5.97%  "

Set-TextValue "D26" 26 "0.175"
$ws.Range("E26").Value = "  -6.71%  "

Set-TextValue "D27" 27 "26.47"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("E28").Value = "  +0.08%  "

Set-TextValue "D29" 29 "7.20"
$ws.Range("E29").Value = "  -4.56%  "

$ws.Range("E30").Value = "  +0.06%  "

Set-TextValue "D31" 31 "6.26"
$ws.Range("E31").Value = "  +2.27%  "

Set-TextValue "D32" 32 "10.03"
$ws.Range("E32").Value = "  -5.59%  "

$ws.Range("E33").Value = "  -5.37%  "

Set-TextValue "D34" 34 "35.46"
$ws.Range("E34").Value = "  -7.73%  "

Set-TextValue "D35" 35 "50.74"
$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D36" 36 "0.0427"
$ws.Range("E36").Value = "  -4.24%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D37" 37 "0.998"
$ws.Range("E37").Value = "  -0.15%  "

Set-TextValue "D38" 38 "3.22"
$ws.Range("E38").Value = "  -1.23%  "

Set-TextValue "D39" 39 "2.83"
$ws.Range("E39").Value = "  +2.98%  "

Set-TextValue "D40" 40 "17.25"
$ws.Range("E40").Value = "  -6.63%  "

$ws.Range("E41").Value = "  -5.69%  "

$ws.Range("E42").Value = "  -4.35%  "

Set-TextValue "D43" 43 "22.84"
$ws.Range("E43").Value = "  -1.01%  "

Set-TextValue "D44" 44 "120.90"
$ws.Range("E44").Value = "  +1.22%  "

Set-TextValue "D45" 45 "2.15"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "2.098.38"
$ws.Range("E46").Value = "  -1.92%  "

Set-TextValue "D47" 47 "3.24"
$ws.Range("E47").Value = "  -7.35%  "

Set-TextValue "D48" 48 "2.30"
$ws.Range("E48").Value = "  -7.35%  "

$ws.Range("D49").Value = "3.231.56"
$ws.Range("E49").Value = "  +0.40%  "

Set-TextValue "D50" 50 "0.238"
$ws.Range("E50").Value = "  -4.99%  "

$ws.Range("E51").Value = "  -5.43%  "
